$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values for rows 2-5 (new dataset)
$data = @{
    2 = @(45127.50694444445, 14.835, 9.791, 3.698, 32.243, 24.166, 11.51, 34.958, 18.033, 7.29, 10.735, 12.533, 13.25, 3.739, 11.655, 16.06, 10.282, 3.096, 1.74, 170.025, 32.298, 10.758, 20.812, 10.713, 2.837, 18.288, 9.502000000000001, 8.641999999999999, 10.303, 12.679, 3.311, 31.418, 5.68, 13.449)
    3 = @(45127.51388888889, 1.868, 0.802, 1.29, 4.347, 2.223, 1.356, 10.693, 2.327, 0.8149999999999999, 0.758, 1.569, 1.696, 0.509, 1.504, 2.024, 1.735, 1.342, 0.595, 15.679, 4.658, 1.388, 2.724, 1.243, 0.823, 5.035, 1.226, 1.328, 1.574, 1.358, 1.246, 10.615, 0.521, 1.744)
    4 = @(45127.52083333334, 13.402, 9.699999999999999, 1.209, 29.415, 23.367, 10.457, 36.895, 16.288, 7.067, 10.36, 11.699, 12.428, 3.378, 10.527, 14.816, 9.083, 0.979, 0.724, 152.872, 29.272, 9.717000000000001, 19.415, 10.213, 1.799, 18.306, 8.583, 7.708, 9.074, 12.093, 0.766, 33.338, 5.35, 12.148)
    5 = @(45127.52777777778, 11.96, 8.73, 0.95, 26.26, 20.98, 9.34, 36.86, 14.54, 6.37, 9.32, 10.46, 11.12, 3.02, 9.4, 13.3, 8.07, 0.75, 0.58, 135.74, 26.31, 8.68, 17.52, 9.17, 1.55, 17.93, 7.66, 6.86, 8.07, 10.84, 0.55, 33.68, 4.8, 10.85)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $col = 1
    foreach ($v in $vals) {
        $ws.Cells.Item([int]$r, $col).Value = $v
        $col = $col + 1
    }
}

# Remove row 6 entirely (dataset now has one fewer row)
$ws.Rows.Item(6).Delete()

# Widen columns whose content now needs more horizontal space
$widthOffset = 0.8333333333333333
$cols8 = @(2,5,6,7,8,9,11,12,13,15,16,17,21,22,23,24,26,29,30,32,34)
foreach ($c in $cols8) {
    $ws.Columns.Item($c).ColumnWidth = 8 - $widthOffset
}
$ws.Columns.Item(20).ColumnWidth = 9 - $widthOffset